$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-02 Sunday" "2024-06-03 Monday"

Replace-Text "185×7=1295" "766×5=3830"
Replace-Text "829×8=6632" "606×7=4242"
Replace-Text "391×4=1564" "894×8=7152"
Replace-Text "973×7=6811" "531×8=4248"
Replace-Text "832×9=7488" "603×8=4824"

Replace-Text "177×7=1239" "111×9=999"
Replace-Text "903×9=8127" "889×2=1778"
Replace-Text "408×2=816" "286×9=2574"
Replace-Text "466×5=2330" "872×2=1744"
Replace-Text "580×2=1160" "514×7=3598"

Replace-Text "982×3=2946" "219×8=1752"
Replace-Text "285×5=1425" "493×3=1479"
Replace-Text "602×7=4214" "310×9=2790"
Replace-Text "927×5=4635" "317×8=2536"
Replace-Text "748×9=6732" "832×3=2496"

Replace-Text "383×2=766" "362×5=1810"
Replace-Text "386×8=3088" "420×2=840"
Replace-Text "563×5=2815" "983×8=7864"
Replace-Text "245×2=490" "944×6=5664"
Replace-Text "792×4=3168" "785×3=2355"

Replace-Text "675×7=4725" "844×7=5908"
Replace-Text "211×7=1477" "323×3=969"
Replace-Text "289×7=2023" "914×8=7312"
Replace-Text "542×2=1084" "683×8=5464"
Replace-Text "374×4=1496" "685×9=6165"
